$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.676.13'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.695.74'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +6.99%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.71'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("E7").Value = '  +7.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.615'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.45%  '

$ws.Range("E10").Value = '  +1.83%  '

$ws.Range("E11").Value = '  +4.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.89'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.02%  '

$ws.Range("E13").Value = '  +1.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.292.77'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +7.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '680.91'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.60%  '

$ws.Range("E16").Value = '  +3.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.695.32'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +7.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '71.789.35'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.28%  '

$ws.Range("E19").Value = '  +2.01%  '

$ws.Range("E20").Value = '  +0.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.63'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.08%  '

$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.946'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.70%  '

$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.36'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +18.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.80'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.95%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '103.30'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.20%  '

$ws.Range("E26").Value = '  +3.04%  '

$ws.Range("E27").Value = '  +4.73%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.32'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.71'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.31'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +5.48%  '

$ws.Range("E31").Value = '  +6.35%  '

$ws.Range("E32").Value = '  +11.15%  '

$ws.Range("E33").Value = '  +1.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '566.63'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("E35").Value = '  +3.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.44'
$ws.Range("D36").ClearFormats()

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.740.05'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.72%  '

$ws.Range("E38").Value = '  -0.09%  '

$ws.Range("E39").Value = '  +2.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0777'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.68'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.46'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.74%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0467'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +9.94%  '

$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.94%  '

$ws.Range("E45").Value = '  +4.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.91'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +8.27%  '

$ws.Range("E47").Value = '  +1.06%  '

$ws.Range("E48").Value = '  +3.13%  '

$ws.Range("E49").Value = '  +1.87%  '

$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '136.49'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.87%  '
